$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.517.51'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').Value = '  -0.57%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.099.71'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').Value = '  +9.76%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '253.01'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +1.26%  '
$ws.Range('E6').Value = '  -6.84%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '47.60'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  +2.81%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '60.20'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  +2.80%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.380'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  +1.97%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0742'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  -2.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.100'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.55'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  -0.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.403.80'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +9.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.833'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  +2.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.100.03'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  +9.73%  '
$ws.Range('E17').Value = '  -0.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.516.38'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.79'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  -2.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0830'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  -3.74%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.21'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -1.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '240.23'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -4.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.24'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  +0.89%  '
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.50'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  -5.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '171.24'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +1.91%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.37'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +13.94%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.15'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +3.99%  '
$ws.Range('E29').Value = '  -9.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '28.82'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +62.78%  '
$ws.Range('E31').Value = '  -5.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.49'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  -1.67%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.990'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  +13.10%  '
$ws.Range('E35').Value = '  +20.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0912'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +2.04%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  -1.14%  '
$ws.Range('E39').Value = '  -5.33%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.35'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -10.68%  '
$ws.Range('E41').Value = '  +6.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0224'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -1.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '98.01'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -6.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.77'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  -5.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.94'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  -9.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.329.25'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -1.55%  '
$ws.Range('E47').Value = '  +3.59%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.10'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +9.90%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.297.35'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  +9.92%  '
$ws.Range('B50').Value = 'MXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.85'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +1.53%  '
$ws.Range('E51').Value = '  -5.82%  '
